$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'60.852.79"
$ws.Cells.Item(2, 4).NumberFormat = "General"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -1.90%  "
$ws.Cells.Item(3, 4).Value = "'3.384.41"
$ws.Cells.Item(3, 4).NumberFormat = "General"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.92%  "
$ws.Cells.Item(5, 4).Value = "'574.76"
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.55%  "
$ws.Cells.Item(6, 4).Value = "'136.65"
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.37%  "
$ws.Cells.Item(7, 5).Value = "  +0.00%  "
$ws.Cells.Item(8, 4).Value = "'3.382.59"
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.95%  "
$ws.Cells.Item(9, 5).Value = "  -1.51%  "
$ws.Cells.Item(10, 4).Value = "'7.57"
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.93%  "
$ws.Cells.Item(11, 5).Value = "  -3.12%  "
$ws.Cells.Item(12, 5).Value = "  -1.44%  "
$ws.Cells.Item(13, 4).Value = "'3.957.83"
$ws.Cells.Item(13, 4).NumberFormat = "General"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.12%  "
$ws.Cells.Item(14, 5).Value = "  +0.85%  "
$ws.Cells.Item(15, 4).Value = "'26.26"
$ws.Cells.Item(15, 4).NumberFormat = "General"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +3.08%  "
$ws.Cells.Item(16, 5).Value = "  -3.39%  "
$ws.Cells.Item(17, 4).Value = "'3.380.99"
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.05%  "
$ws.Cells.Item(18, 4).Value = "'60.934.72"
$ws.Cells.Item(18, 4).NumberFormat = "General"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.79%  "
$ws.Cells.Item(19, 4).Value = "'14.05"
$ws.Cells.Item(19, 4).NumberFormat = "General"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.69%  "
$ws.Cells.Item(20, 5).Value = "  -0.94%  "
$ws.Cells.Item(21, 4).Value = "'9.49"
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.23%  "
$ws.Cells.Item(22, 4).Value = "'377.79"
$ws.Cells.Item(22, 4).NumberFormat = "General"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.40%  "
$ws.Cells.Item(23, 5).Value = "  -2.72%  "
$ws.Cells.Item(24, 4).Value = "'3.526.44"
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.77%  "
$ws.Cells.Item(25, 5).Value = "  -0.04%  "
$ws.Cells.Item(26, 5).Value = "  -2.42%  "
$ws.Cells.Item(27, 4).Value = "'71.42"
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.26%  "
$ws.Cells.Item(28, 5).Value = "  +11.09%  "
$ws.Cells.Item(29, 5).Value = "  +6.11%  "
$ws.Cells.Item(30, 4).Value = "'7.55"
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.13%  "
$ws.Cells.Item(31, 5).Value = "  +0.09%  "
$ws.Cells.Item(32, 4).Value = "'8.18"
$ws.Cells.Item(32, 4).NumberFormat = "General"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.65%  "
$ws.Cells.Item(33, 5).Value = "  -0.70%  "
$ws.Cells.Item(34, 5).Value = "  +0.03%  "
$ws.Cells.Item(35, 4).Value = "'23.75"
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.72%  "
$ws.Cells.Item(36, 4).Value = "'5.22"
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -4.75%  "
$ws.Cells.Item(37, 4).Value = "'6.86"
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.15%  "
$ws.Cells.Item(38, 4).Value = "'1.54"
$ws.Cells.Item(38, 4).NumberFormat = "General"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.40%  "
$ws.Cells.Item(39, 4).Value = "'164.61"
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.75%  "
$ws.Cells.Item(40, 4).Value = "'0.0758"
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -4.32%  "
$ws.Cells.Item(41, 4).Value = "'0.999"
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.05%  "
$ws.Cells.Item(42, 4).Value = "'0.772"
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.38%  "
$ws.Cells.Item(43, 5).Value = "  -3.13%  "
$ws.Cells.Item(44, 5).Value = "  -1.19%  "
$ws.Cells.Item(45, 4).Value = "'41.60"
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.44%  "
$ws.Cells.Item(46, 5).Value = "  -2.65%  "
$ws.Cells.Item(47, 4).Value = "'24.09"
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -3.56%  "
$ws.Cells.Item(48, 4).Value = "'23.44"
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.37%  "
$ws.Cells.Item(49, 4).Value = "'2.439.31"
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.32%  "
$ws.Cells.Item(50, 5).Value = "  -2.35%  "
$ws.Cells.Item(51, 4).Value = "'2.39"
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +4.54%  "
